$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 929.7059
$ws.Range("I28").Value = 913.26666
$ws.Range("J28").Value = 1053
$ws.Range("K28").Value = 913.26666
$ws.Range("L28").Value = 1053
$ws.Range("M28").Value = -428.26666
$ws.Range("N28").Value = -2023

$ws.Range("H43").Value = 1910
$ws.Range("I43").Value = 1776
$ws.Range("J43").Value = 2133.3333
$ws.Range("K43").Value = 1776
$ws.Range("L43").Value = 2133.3333
$ws.Range("M43").Value = -1707
$ws.Range("N43").Value = -2271.3333

$ws.Range("H63").Value = 199900
$ws.Range("J63").Value = 199900
$ws.Range("L63").Value = 199900
$ws.Range("N63").Value = -201148

$ws.Range("H66").Value = 199900
$ws.Range("J66").Value = 199900
$ws.Range("L66").Value = 599700
$ws.Range("N66").Value = -605940

$ws.Range("N107").ClearContents()
$ws.Range("H107").Value = 924.1875
$ws.Range("I107").Value = 924.1875
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 924.1875
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 995.8125

$ws.Range("H113").Value = 13778.111
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 14875.375
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 14875.375
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = -21383.375

$ws.Range("H134").Value = 52499.89
$ws.Range("J134").Value = 52499.89
$ws.Range("L134").Value = 52499.89
$ws.Range("N134").Value = -62639.89

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11604.515
$ws.Range("I32").Value = 7653.9
$ws.Range("J32").Value = 16872
$ws.Range("K32").Value = 7653.9
$ws.Range("L32").Value = 16872
$ws.Range("M32").Value = -7366.9
$ws.Range("N32").Value = -17446

$ws.Range("N110").ClearContents()
$ws.Range("H110").Value = 3543.1428
$ws.Range("I110").Value = 3543.1428
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 3543.1428
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1498.1428

$ws.Range("H122").Value = 3063.7144
$ws.Range("I122").Value = 1768
$ws.Range("J122").Value = 5396
$ws.Range("K122").Value = 5304
$ws.Range("L122").Value = 16188
$ws.Range("M122").Value = -2854
$ws.Range("N122").Value = -21088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6971.346
$ws.Range("I20").Value = 1435.375
$ws.Range("J20").Value = 15828.9
$ws.Range("K20").Value = 1435.375
$ws.Range("L20").Value = 15828.9
$ws.Range("M20").Value = -1188.375
$ws.Range("N20").Value = -16322.9

$ws.Range("H47").Value = 79800
$ws.Range("J47").Value = 79800
$ws.Range("L47").Value = 79800
$ws.Range("N47").Value = -80840

$ws.Range("H107").Value = 1400
$ws.Range("I107").Value = 1400
$ws.Range("K107").Value = 1400
$ws.Range("M107").Value = 520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5850079.5
$ws.Range("I16").Value = 18519534
$ws.Range("J16").Value = 2639.3076
$ws.Range("K16").Value = 18519534
$ws.Range("L16").Value = 2639.3076
$ws.Range("M16").Value = -18519247
$ws.Range("N16").Value = -3213.3076

$ws.Range("H31").Value = 3316.8333
$ws.Range("I31").Value = 1163.5
$ws.Range("J31").Value = 6546.8335
$ws.Range("K31").Value = 1163.5
$ws.Range("L31").Value = 6546.8335
$ws.Range("M31").Value = -868.5
$ws.Range("N31").Value = -7136.8335

$ws.Range("H34").Value = 3316.8333
$ws.Range("I34").Value = 1163.5
$ws.Range("J34").Value = 6546.8335
$ws.Range("K34").Value = 1163.5
$ws.Range("L34").Value = 6546.8335
$ws.Range("M34").Value = -961.5
$ws.Range("N34").Value = -6950.8335

$ws.Range("H99").Value = 4535.077
$ws.Range("I99").Value = 2549.5557
$ws.Range("J99").Value = 9002.5
$ws.Range("K99").Value = 2549.5557
$ws.Range("L99").Value = 9002.5
$ws.Range("M99").Value = -1051.5557
$ws.Range("N99").Value = -11998.5

$ws.Range("H113").Value = 5850079.5
$ws.Range("I113").Value = 18519534
$ws.Range("J113").Value = 2639.3076
$ws.Range("K113").Value = 18519534
$ws.Range("L113").Value = 2639.3076
$ws.Range("M113").Value = -18517364
$ws.Range("N113").Value = -6979.3076

$ws.Range("H126").Value = 4535.077
$ws.Range("I126").Value = 2549.5557
$ws.Range("J126").Value = 9002.5
$ws.Range("K126").Value = 7648.6671
$ws.Range("L126").Value = 27007.5
$ws.Range("M126").Value = -5178.6671
$ws.Range("N126").Value = -31947.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2322.92
$ws.Range("I129").Value = 4732.5
$ws.Range("J129").Value = 1189
$ws.Range("K129").Value = 14197.5
$ws.Range("L129").Value = 3567
$ws.Range("M129").Value = -9197.5
$ws.Range("N129").Value = -13567

$ws.Range("H132").Value = 2570.6538
$ws.Range("I132").Value = 1065.4445
$ws.Range("J132").Value = 3367.5293
$ws.Range("K132").Value = 9589.0005
$ws.Range("L132").Value = 30307.7637
$ws.Range("M132").Value = -7059.0005
$ws.Range("N132").Value = -35367.7637

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 15873829
$ws.Range("I107").Value = 700.5
$ws.Range("J107").Value = 37038000
$ws.Range("K107").Value = 700.5
$ws.Range("L107").Value = 37038000
$ws.Range("M107").Value = 1219.5
$ws.Range("N107").Value = -37041840

$ws.Range("H122").Value = 4222.737
$ws.Range("I122").Value = 3625.4167
$ws.Range("J122").Value = 5246.7144
$ws.Range("K122").Value = 10876.2501
$ws.Range("L122").Value = 15740.1432
$ws.Range("M122").Value = -8426.250100000001
$ws.Range("N122").Value = -20640.1432

$ws.Range("H123").Value = 15261.277
$ws.Range("J123").Value = 15261.277
$ws.Range("L123").Value = 15261.277
$ws.Range("N123").Value = -20161.277

$ws.Range("H136").Value = 26319.723
$ws.Range("J136").Value = 26319.723
$ws.Range("L136").Value = 78959.16900000001
$ws.Range("N136").Value = -84059.16900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 19900
$ws.Range("J64").Value = 19900
$ws.Range("L64").Value = 19900
$ws.Range("N64").Value = -20396

$ws.Range("H67").Value = 19900
$ws.Range("J67").Value = 19900
$ws.Range("L67").Value = 19900
$ws.Range("N67").Value = -21616

$ws.Range("H113").Value = 14672.286
$ws.Range("I113").Value = 25338.5
$ws.Range("J113").Value = 450.66666
$ws.Range("K113").Value = 76015.5
$ws.Range("L113").Value = 1351.99998
$ws.Range("M113").Value = -73845.5
$ws.Range("N113").Value = -5691.999980000001
